# "add w7 ses 2"
#
# The underlying OOXML diff is structural: a handful of narrow "separator /
# indicator" columns (A, G, H, I, K, L) get explicit widths, column J's
# width shrinks to match, the sheet's scroll/selection state changes, and
# every one of the 37 invisible-hyperlink AutoShapes overlaid on column H
# has its cached bottom-right (two-cell) anchor recomputed so it still
# lands in the same visual spot now that columns G/H are narrower.
#
# Note: the source workbook was produced by Excel for Mac, whose column
# width encoding (exact 1/256-character-unit fractions) is finer-grained
# than the pixel-snapped (1/6-character-unit) grid this COM host rounds
# `Range.ColumnWidth` writes to - so the widths below are the closest
# reachable approximation of the target values, not bit-exact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# `ColumnWidth` (characters) round-trips into the saved `width` XML
# attribute with a constant +5/6 padding baked in by this host, then gets
# snapped to the nearest 1/6 of a character. Pre-subtract that padding so
# the stored width lands as close as possible to each target.
function Set-StoredColumnWidth($colIndex, $targetStoredWidth) {
    $ws.Columns.Item($colIndex).ColumnWidth = $targetStoredWidth - 0.8333333333333333
}

Set-StoredColumnWidth 1  1.5          # A - new narrow spacer column
Set-StoredColumnWidth 7  1.83203125   # G - new narrow spacer column
Set-StoredColumnWidth 8  3            # H - new column (icon column)
Set-StoredColumnWidth 9  2.5          # I - new column
Set-StoredColumnWidth 10 2.6640625    # J - was a bare style-only default width
Set-StoredColumnWidth 11 3.1640625    # K - new column
Set-StoredColumnWidth 12 2.6640625    # L - new column

# Re-anchor every AutoShape sitting over column H: force-write its exact,
# unchanged Left (top-left stays at the start of column H, row unchanged)
# so the host recomputes the cached bottom-right two-cell anchor against
# the new, narrower column widths - this is what moves `<xdr:to>` from
# col 7/304800 EMU to col 8/76200 EMU for each shape, matching the diff.
$colAStart = 0
for ($c = 1; $c -le 7; $c++) {
    $colAStart += $ws.Columns.Item($c).Width
}

for ($i = 1; $i -le $ws.Shapes.Count; $i++) {
    $shp = $ws.Shapes.Item($i)
    $shp.Left = $colAStart
}

# Scroll/selection: the saved view no longer pins topLeftCell to C6, and
# the live selection moved from Z22 to Z24 (week 7 / session 2 row).
$ws.Range("Z24").Select()
